$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-30 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-01 Monday", 2) | Out-Null
$d.Content.Find.Execute("664×9=5976", $true, $false, $false, $false, $false, $true, 1, $false, "412×7=2884", 2) | Out-Null
$d.Content.Find.Execute("284×4=1136", $true, $false, $false, $false, $false, $true, 1, $false, "830×8=6640", 2) | Out-Null
$d.Content.Find.Execute("934×2=1868", $true, $false, $false, $false, $false, $true, 1, $false, "260×8=2080", 2) | Out-Null
$d.Content.Find.Execute("704×7=4928", $true, $false, $false, $false, $false, $true, 1, $false, "524×7=3668", 2) | Out-Null
$d.Content.Find.Execute("965×8=7720", $true, $false, $false, $false, $false, $true, 1, $false, "489×5=2445", 2) | Out-Null
$d.Content.Find.Execute("919×7=6433", $true, $false, $false, $false, $false, $true, 1, $false, "471×8=3768", 2) | Out-Null
$d.Content.Find.Execute("708×2=1416", $true, $false, $false, $false, $false, $true, 1, $false, "798×3=2394", 2) | Out-Null
$d.Content.Find.Execute("262×9=2358", $true, $false, $false, $false, $false, $true, 1, $false, "953×9=8577", 2) | Out-Null
$d.Content.Find.Execute("723×9=6507", $true, $false, $false, $false, $false, $true, 1, $false, "900×2=1800", 2) | Out-Null
$d.Content.Find.Execute("936×5=4680", $true, $false, $false, $false, $false, $true, 1, $false, "435×7=3045", 2) | Out-Null
$d.Content.Find.Execute("161×9=1449", $true, $false, $false, $false, $false, $true, 1, $false, "855×3=2565", 2) | Out-Null
$d.Content.Find.Execute("821×3=2463", $true, $false, $false, $false, $false, $true, 1, $false, "241×3=723", 2) | Out-Null
$d.Content.Find.Execute("886×6=5316", $true, $false, $false, $false, $false, $true, 1, $false, "927×7=6489", 2) | Out-Null
$d.Content.Find.Execute("498×9=4482", $true, $false, $false, $false, $false, $true, 1, $false, "826×6=4956", 2) | Out-Null
$d.Content.Find.Execute("772×4=3088", $true, $false, $false, $false, $false, $true, 1, $false, "206×3=618", 2) | Out-Null
$d.Content.Find.Execute("761×9=6849", $true, $false, $false, $false, $false, $true, 1, $false, "355×5=1775", 2) | Out-Null
$d.Content.Find.Execute("504×4=2016", $true, $false, $false, $false, $false, $true, 1, $false, "126×3=378", 2) | Out-Null
$d.Content.Find.Execute("176×7=1232", $true, $false, $false, $false, $false, $true, 1, $false, "823×5=4115", 2) | Out-Null
$d.Content.Find.Execute("472×3=1416", $true, $false, $false, $false, $false, $true, 1, $false, "378×6=2268", 2) | Out-Null
$d.Content.Find.Execute("844×6=5064", $true, $false, $false, $false, $false, $true, 1, $false, "519×5=2595", 2) | Out-Null
$d.Content.Find.Execute("690×9=6210", $true, $false, $false, $false, $false, $true, 1, $false, "289×6=1734", 2) | Out-Null
$d.Content.Find.Execute("124×8=992", $true, $false, $false, $false, $false, $true, 1, $false, "567×6=3402", 2) | Out-Null
$d.Content.Find.Execute("223×3=669", $true, $false, $false, $false, $false, $true, 1, $false, "446×8=3568", 2) | Out-Null
$d.Content.Find.Execute("291×6=1746", $true, $false, $false, $false, $false, $true, 1, $false, "517×5=2585", 2) | Out-Null
$d.Content.Find.Execute("830×4=3320", $true, $false, $false, $false, $false, $true, 1, $false, "354×3=1062", 2) | Out-Null
